$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Add a "Date: " label in front of the {{Form3NDate}} merge field that
#    precedes the "Signed:" line. The new run is plain (not bold), using
#    the same Times New Roman / 14pt / en-US formatting as the rest of
#    the document, while the existing bold "{{Form3NDate}}" field itself
#    is left untouched.
# -----------------------------------------------------------------------
$dateField = $d.Content
$dateField.Find.Execute("{{Form3NDate}}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)

$dateLabel = $dateField.Duplicate
$dateLabel.Collapse(1)                 # wdCollapseStart
$dateLabel.InsertBefore("Date: ")

$dateLabel.Font.Name = "Times New Roman"
$dateLabel.Font.Size = 14
$dateLabel.Font.Bold = $false
$dateLabel.LanguageID = "en-US"

# -----------------------------------------------------------------------
# 2) Turn the "Signed: {{Form3MyName}}" merge field into a blank
#    signature line made of underscores (digital-signing field removed).
# -----------------------------------------------------------------------
$signatureField = $d.Content
$signatureField.Find.Execute(" {{Form3MyName}}", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
$signatureField.Text = "_________________________"
